$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: "Productivity" -> "Plan Vs. Actual", with new raw numbers and flipped formula
$ws.Range("A2").Value = "Plan Vs. Actual"
$ws.Range("B2").Value = 8000
$ws.Range("C2").Value = 10000
$ws.Range("D2").Formula = "=B2-C2"

# Row 3: "Efficiency" values become percentages, formula sign flips
$ws.Range("B3").Value = 0.65
$ws.Range("C3").Value = 0.7
$ws.Range("D3").Formula = "=B3-C3"
$ws.Range("B3:D3").NumberFormat = "0%"

# Row 4: "Lost Time" values become percentages, shared formula follows row 3's pattern
$ws.Range("B4").Value = 0.04
$ws.Range("C4").Value = 0.03
$ws.Range("D4").Formula = "=B4-C4"
$ws.Range("B4:D4").NumberFormat = "0%"

# Update the active selection shown when the workbook is reopened
$ws.Range("G10").Select()
